# Update the waves sheet so the enemy lookup column stores numeric
# template ids (tid) instead of the old "enemy:slug" string keys.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: F4 "string" -> "tid"
$ws.Range("F4").Value = "tid"

# Force column F data rows to stay text-formatted (the sheet stores all
# values, including numeric-looking ones, as text - see the
# numberStoredAsText ignored-error range over A4:J11).
$ws.Range("F6:F11").NumberFormat = "@"

$ws.Range("F6").Value  = "40060001"
$ws.Range("F7").Value  = "40060002"
$ws.Range("F8").Value  = "40060003"
$ws.Range("F9").Value  = "40060004"
$ws.Range("F10").Value = "40060006"
$ws.Range("F11").Value = "40060005"

Write-Host "Applied enemy tid lookup changes"
